{"js": "// 1) Remove the whole bullet paragraph \"Serie de \u00edndice de masculinidad por\n//    grupo quinquenal etario\" (under \"Mujeres y poblaci\u00f3n\").\nconst serieResults = context.document.body.search(\n  \"Serie de \u00edndice de masculinidad por grupo quinquenal etario\",\n  { matchCase: true }\n);\nserieResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < serieResults.items.length; i++) {\n  serieResults.items[i].paragraphs.getFirst().delete();\n}\nawait context.sync();\n\n// 2) Append \" (brecha)\" right after \"Tasa de analfabetismo por sexo\".\nconst analfabetismoResults = context.document.body.search(\n  \"Tasa de analfabetismo por sexo\",\n  { matchCase: true }\n);\nanalfabetismoResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < analfabetismoResults.items.length; i++) {\n  analfabetismoResults.items[i].insertText(\" (brecha)\", Word.InsertLocation.after);\n}\nawait context.sync();\n\n// 3) Remove the two bullet paragraphs under \"Mujeres y sociedad\":\n//    \"Proporci\u00f3n de mujeres en edad f\u00e9rtil (Mapa departamental)\" and\n//    \"Distribuci\u00f3n de mujeres en edad f\u00e9rtil por grupo quinquenal etario\".\nconst proporcionResults = context.document.body.search(\n  \"Proporci\u00f3n de mujeres en edad f\u00e9rtil (Mapa departamental)\",\n  { matchCase: true }\n);\nproporcionResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < proporcionResults.items.length; i++) {\n  proporcionResults.items[i].paragraphs.getFirst().delete();\n}\nawait context.sync();\n\nconst distribucionResults = context.document.body.search(\n  \"Distribuci\u00f3n de mujeres en edad f\u00e9rtil por grupo quinquenal etario\",\n  { matchCase: true }\n);\ndistribucionResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < distribucionResults.items.length; i++) {\n  distribucionResults.items[i].paragraphs.getFirst().delete();\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Remove the whole bullet paragraph \"Serie de \u00edndice de masculinidad por\n#    grupo quinquenal etario\" (under \"Mujeres y poblaci\u00f3n\").\n$range = $d.Content\nif ($range.Find.Execute(\"Serie de \u00edndice de masculinidad por grupo quinquenal etario\")) {\n    $range.Paragraphs(1).Range.Delete()\n}\n\n# 2) Append \" (brecha)\" right after \"Tasa de analfabetismo por sexo\".\n$range = $d.Content\nif ($range.Find.Execute(\"Tasa de analfabetismo por sexo\")) {\n    $range.Collapse(0)\n    $range.InsertAfter(\" (brecha)\")\n}\n\n# 3) Remove the two bullet paragraphs under \"Mujeres y sociedad\":\n#    \"Proporci\u00f3n de mujeres en edad f\u00e9rtil (Mapa departamental)\" and\n#    \"Distribuci\u00f3n de mujeres en edad f\u00e9rtil por grupo quinquenal etario\".\n$range = $d.Content\nif ($range.Find.Execute(\"Proporci\u00f3n de mujeres en edad f\u00e9rtil (Mapa departamental)\")) {\n    $range.Paragraphs(1).Range.Delete()\n}\n\n$range = $d.Content\nif ($range.Find.Execute(\"Distribuci\u00f3n de mujeres en edad f\u00e9rtil por grupo quinquenal etario\")) {\n    $range.Paragraphs(1).Range.Delete()\n}\n"}
